# Inserts a new weekly price record for "Pepino dulce" (Macroferia Regional
# de Talca) at row 24 — pushing the existing rows 24-46 down to 25-47 — and
# populates the new row with the latest sampled values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 24 (shifts rows 24:46 -> 25:47).
$ws.Rows.Item(24).EntireRow.Insert()

# Fill the newly inserted row 24 with the market record.
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 44791
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 100112043
$ws.Range("G24").Value = "Pepino dulce"
$ws.Range("H24").Value = "Cultivar IV Región"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 15000
$ws.Range("N24").Value = "$/bandeja 18 kilos"
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 833
$ws.Range("Q24").Value = 18
$ws.Range("R24").Value = "Hortaliza"
